# The workbook's sector/group codelist has its "category" and "group"
# columns out of order. This swaps column D <-> E (category-name /
# group-name) and column F <-> G (category-code / group-code) for every
# row, including the header row, using a scratch helper column (I) and
# PasteSpecial (values only) so that text-typed cells (e.g. numeric-looking
# codes like "110") keep their original string type instead of being
# coerced to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$xlPasteValues = -4163

$helperCol = "I"
$helperRange = $helperCol + "1:" + $helperCol + $lastRow

# Swap D <-> E
$ws.Range("D1:D$lastRow").Copy()
$ws.Range($helperCol + "1").PasteSpecial($xlPasteValues)
$ws.Range("E1:E$lastRow").Copy()
$ws.Range("D1").PasteSpecial($xlPasteValues)
$ws.Range($helperRange).Copy()
$ws.Range("E1").PasteSpecial($xlPasteValues)

# Swap F <-> G
$ws.Range("F1:F$lastRow").Copy()
$ws.Range($helperCol + "1").PasteSpecial($xlPasteValues)
$ws.Range("G1:G$lastRow").Copy()
$ws.Range("F1").PasteSpecial($xlPasteValues)
$ws.Range($helperRange).Copy()
$ws.Range("G1").PasteSpecial($xlPasteValues)

# Clean up helper column
$ws.Range($helperRange).Clear()

$excel.CutCopyMode = $false
